# Update "想去人数" (want-to-go count) values in the F column of the
# "展览" and "全部类型" sheets, as per refreshed output data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value, identical on both sheets.
$updates = @{
    2  = 223
    3  = 438
    4  = 12979
    5  = 1338
    16 = 46
    18 = 5526
    24 = 137
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
